$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 314 (pushes existing rows 314:441 down to 316:443)
$ws.Rows("314:315").Insert()

# New row 314 (Primera) - new weekly report, date 2022-07-13 (serial 44755)
$ws.Range("A314").Value = 8
$ws.Range("B314").Value = "Terminal La Palmera de La Serena"
$ws.Range("C314").Value = "Coquimbo"
$ws.Range("D314").Value = 44755
$ws.Range("E314").Value = 4
$ws.Range("F314").Value = 100112009
$ws.Range("G314").Value = "Acelga"
$ws.Range("H314").Value = "Sin especificar"
$ws.Range("I314").Value = "Primera"
$ws.Range("J314").Value = 2520
$ws.Range("K314").Value = 600
$ws.Range("L314").Value = 700
$ws.Range("M314").Value = 650
$ws.Range("N314").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O314").Value = "Provincia del Elquí"
$ws.Range("P314").Value = 325
$ws.Range("Q314").Value = 2
$ws.Range("R314").Value = "Hortaliza"

# New row 315 (Segunda) - same date
$ws.Range("A315").Value = 8
$ws.Range("B315").Value = "Terminal La Palmera de La Serena"
$ws.Range("C315").Value = "Coquimbo"
$ws.Range("D315").Value = 44755
$ws.Range("E315").Value = 4
$ws.Range("F315").Value = 100112009
$ws.Range("G315").Value = "Acelga"
$ws.Range("H315").Value = "Sin especificar"
$ws.Range("I315").Value = "Segunda"
$ws.Range("J315").Value = 1340
$ws.Range("K315").Value = 500
$ws.Range("L315").Value = 550
$ws.Range("M315").Value = 525
$ws.Range("N315").Value = "`$/atado 1,5 a 2 kilos"
$ws.Range("O315").Value = "Provincia del Elquí"
$ws.Range("P315").Value = 262
$ws.Range("Q315").Value = 2
$ws.Range("R315").Value = "Hortaliza"
